$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update R column values (50 -> 20) for these "line" rows ---
$rowsToTwenty = @(11, 12, 29, 30, 31, 32, 33, 43, 54, 55)
foreach ($r in $rowsToTwenty) {
    $ws.Range("R$r").Value = 20
}

# --- Remove Q/R values entirely for rows 50, 58, 59 ---
$rowsToClear = @(50, 58, 59)
foreach ($r in $rowsToClear) {
    $ws.Range("Q$r").ClearContents()
    $ws.Range("R$r").ClearContents()
}

# --- New column width for column F (enum line names column) ---
$ws.Columns.Item(6).ColumnWidth = 22.05

# --- View: zoom, scroll position, and active selection ---
$win = $excel.ActiveWindow
$win.Zoom = 70
$win.ScrollRow = 28
$win.ScrollColumn = 1
$ws.Range("B62").Select()
